$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# New case row (id 43) appended as row 44, mirroring the layout of the
# existing rows (id/type/added/SocialMedia/source/motivation/lang/lookalike/description/Entity).
$row = 44

$ws.Cells.Item($row, 1).Value = 43
$ws.Cells.Item($row, 2).Value = "msg"

# Match the date style used by the other rows (style index carrying the
# mm/dd/yyyy-style date format) by copying formatting from C43 before
# writing the date value.
$ws.Cells.Item(43, 3).Copy()
$ws.Cells.Item($row, 3).PasteSpecial(-4122)
$ws.Cells.Item($row, 3).Value = Get-Date -Year 2021 -Month 7 -Day 7 -Hour 0 -Minute 0 -Second 0

$ws.Cells.Item($row, 4).Value = "SocialMedia"
$ws.Cells.Item($row, 5).Value = "unknown"
$ws.Cells.Item($row, 6).Value = "opportunity"
$ws.Cells.Item($row, 7).Value = "en"
$ws.Cells.Item($row, 8).Value = "no"
$ws.Cells.Item($row, 9).Value = "busines opportunity"
$ws.Cells.Item($row, 10).Value = "Lombard"

$ws.Range("F44").Select()
